$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "46.956.35"
$ws.Range("E2").Value = "  +5.44%  "
$ws.Range("D3").Value = "2.342.67"
$ws.Range("E3").Value = "  +4.33%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.85%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.62"
$ws.Range("E5").Value = "  +0.48%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "99.17"
$ws.Range("E6").Value = "  +4.91%  "
$ws.Range("E7").Value = "  +1.83%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  -0.59%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.538"
$ws.Range("E9").Value = "  +4.19%  "
$ws.Range("E10").Value = "  +4.09%  "
$ws.Range("E11").Value = "  +1.02%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.44"
$ws.Range("E12").Value = "  +3.62%  "
$ws.Range("E13").Value = "  -0.16%  "
$ws.Range("D14").Value = "2.696.55"
$ws.Range("E14").Value = "  +4.15%  "
$ws.Range("D15").Value = "2.337.18"
$ws.Range("E15").Value = "  +4.33%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.21"
$ws.Range("E16").Value = "  +4.98%  "
$ws.Range("E17").Value = "  +0.00%  "
$ws.Range("D18").Value = "46.855.06"
$ws.Range("E18").Value = "  +5.68%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.36"
$ws.Range("E19").Value = "  +14.20%  "
$ws.Range("E20").Value = "  +1.38%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.18"
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "66.76"
$ws.Range("E22").Value = "  +2.44%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "245.83"
$ws.Range("E23").Value = "  +3.79%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.98"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.97"
$ws.Range("E25").Value = "  +0.46%  "
$ws.Range("E26").Value = "  -0.37%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "42.05"
$ws.Range("E27").Value = "  +13.50%  "
$ws.Range("E28").Value = "  -0.50%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.88"
$ws.Range("E29").Value = "  +1.48%  "
$ws.Range("E30").Value = "  +1.50%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.76"
$ws.Range("E31").Value = "  -2.17%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "151.28"
$ws.Range("E32").Value = "  +1.62%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0811"
$ws.Range("E33").Value = "  +3.81%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.62"
$ws.Range("E34").Value = "  +0.75%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.05"
$ws.Range("E35").Value = "  -4.96%  "
$ws.Range("E36").Value = "  +0.74%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.108"
$ws.Range("E37").Value = "  -0.69%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.82"
$ws.Range("E38").Value = "  -2.21%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.03"
$ws.Range("E39").Value = "  +6.70%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0319"
$ws.Range("E40").Value = "  +7.16%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.45"
$ws.Range("E41").Value = "  +3.05%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "13.86"
$ws.Range("E42").Value = "  -8.96%  "
$ws.Range("E43").Value = "  -0.71%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.98"
$ws.Range("E44").Value = "  +12.78%  "
$ws.Range("D45").Value = "1.810.56"
$ws.Range("E45").Value = "  +0.33%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.198"
$ws.Range("E46").Value = "  +6.33%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "81.21"
$ws.Range("E47").Value = "  -0.71%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "73.56"
$ws.Range("E48").Value = "  +7.40%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.93"
$ws.Range("E49").Value = "  +2.07%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "98.37"
$ws.Range("E50").Value = "  +0.10%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "55.55"
$ws.Range("E51").Value = "  +3.50%  "
